$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set demscalar_ippu_* rows (4-10,12-15,21-22) columns J:AS (col 10 to col 45) to 1 ---
$onesRows = @(4,5,6,7,8,9,10,12,13,14,15,21,22)
foreach ($r in $onesRows) {
    $ws.Range($ws.Cells.Item($r,10), $ws.Cells.Item($r,45)).Value = 1
}

# --- Set prodinit_ippu_* rows (96-99,101-104,111-112) columns J:AS (col 10 to col 45) to new calculated values ---
# row 96
$vals96 = @(2957017.89188694,3197087.51714342,3195934.89355052,2540308.18177388,2282861.43504041,2539715.41786911,2539715,2716635.51723646,2905880.59428346,3108308.77924437,3324838.42802565,3556451.87064183,3804199.86775195,4069206.37764521,4352673.6553023,4655887.70666538,4980224.12286095,5327154.32084381,5698252.21877414,6095201.37641164,6519802.63292079,6973982.2767375,7459800.78456153,7979462.16912065,8535323.97811481,9129907.98970226,9765911.65305008,10446220.3258512,11173920.3643263,11952313.1250949,12784929.9424402,13675548.1489115,14628208.2119491,15647232.0642714,16737242.7111864,17903185.2037792)
for ($i = 0; $i -lt $vals96.Length; $i++) {
    $ws.Cells.Item(96, 10+$i).Value = $vals96[$i]
}

# row 97
$vals97 = @(27862.1278568579,24788.9903959027,30942.558746724,23404.9614863768,24454.2738572532,30290.6533790425,30290.65,32428.1027170808,34716.384291177,37166.1379257276,39788.757859362,42596.442362526,45602.2504737998,48820.1627398034,52265.146241198,55953.2242072328,59901.5505426261,64128.4896134182,68653.7016628951,73498.2342548668,78684.6201696196,84236.982207873,90181.145390205,96544.7570738021,103357.415545222,110650.807687271,118458.856360322,126817.878183547,135766.752449957,145347.101960892,155603.486621047,166583.610694495,178338.54468568,190922.962877387,204395.397630533,218818.511628555)
for ($i = 0; $i -lt $vals97.Length; $i++) {
    $ws.Cells.Item(97, 10+$i).Value = $vals97[$i]
}

# row 98
$vals98 = @(99441.9885513451,116221.061294563,77581.1000209051,66776.3942495914,61425.0077572945,65645.8378229215,65645.84,71277.725016863,77392.7804652898,84032.4584816857,91241.7674623183,99069.5771618077,107568.949966612,116797.500589114,126817.786620565,137697.732590334,149511.090406024,162337.939295649,176265.228640832,191387.367380717,207806.863981999,225635.021313241,244992.691133805,266011.093311878,288832.705324803,313612.228071389,340517.634543152,369731.308463102,401451.280610591,435892.571212878,473288.6475031,513893.006324986,557980.892512327,605851.164691585,657828.321155354,714264.699539468)
for ($i = 0; $i -lt $vals98.Length; $i++) {
    $ws.Cells.Item(98, 10+$i).Value = $vals98[$i]
}

# row 99
$vals99 = @(2996.4025639497,3129.75344093452,3445.59211543944,3145.17133732205,2648.66151780243,2570.04318161036,2570.043,2687.79935431454,2810.95116659669,2939.74565039896,3074.44134630871,3215.30864093966,3362.63030970282,3516.70208444674,3677.83324710619,3846.34725055128,4022.58236788304,4206.89237147903,4399.64724315201,4601.23391684719,4812.05705536913,5032.53986269719,5263.12493352044,5504.27514169724,5756.47456942293,6020.22947897098,6296.06932895813,6584.54783717364,6886.24409210605,7201.7637153989,7531.74007756876,7876.83556942598,8237.74293175053,8615.18664589201,9009.92438808525,9422.74855040105)
for ($i = 0; $i -lt $vals99.Length; $i++) {
    $ws.Cells.Item(99, 10+$i).Value = $vals99[$i]
}

# row 101
$vals101 = @(168757.46913669,189820.508197794,185606.982195615,171058.604878699,166840.211783412,188949.060947327,188949.1,199872.081328056,211426.51060211,223648.891272687,236577.837022691,250254.193758081,264721.168650723,280024.466641098,296212.434832106,313336.215230156,331449.906316078,350610.733956307,370879.232194303,392319.434493347,414999.076034929,438989.80771182,464367.422491903,491212.094867885,519608.634149402,549646.752397707,581421.347849429,615032.804724814,650587.31036761,688197.19071854,727981.265182189,770065.222008444,814582.015374403,861672.285421231,911484.802572992,964176.937541149)
for ($i = 0; $i -lt $vals101.Length; $i++) {
    $ws.Cells.Item(101, 10+$i).Value = $vals101[$i]
}

# row 102
$vals102 = @(17383.7934995589,14935.6222411441,39203.4152876161,31047.5146410628,4920.73668115729,2359.88049545489,2359.88,5432.43566584759,12505.4482700701,28787.4990252643,66268.7240179154,152550.375394145,351170.440924789,808393.150529368,1860918.25981094,4283827.45132124,9861356.09984172,22700807.9183704,52257181.9665875,120295853.650176,276920642.499957,637470369.224788,1467454603.49584,3378075464.02489,7776318131.72426,17901057667.2357,41208173351.1388,94861073714.4079,218369866326.57,502686683296.96,1157183020784.55,2663831344823.98,6132130619109.68,14116143652595.7,32495314271314.1,74804101996886.7)
for ($i = 0; $i -lt $vals102.Length; $i++) {
    $ws.Cells.Item(102, 10+$i).Value = $vals102[$i]
}

# row 103
$vals103 = @(114561.556166575,115679.940240887,132774.810717697,119436.604949319,158470.459498445,117278.899584277,117278.9,126167.328269109,135729.399939504,146016.169642933,157082.561380929,168987.661779742,181795.035570711,195573.064980523,210395.314843449,226340.92538501,243495.034774238,261949.233700686,281802.054403273,303159.496762052,326135.594261847,350853.022849595,377443.755936269,406049.769040609,436823.797836972,469930.153654717,505545.600783325,543860.300267428,585078.825230939,629421.25315139,677124.340916442,728442.788936534,783650.601063115,843042.547573473,906935.739033435,975671.319441249)
for ($i = 0; $i -lt $vals103.Length; $i++) {
    $ws.Cells.Item(103, 10+$i).Value = $vals103[$i]
}

# row 104
$vals104 = @(42332.9293197324,48734.0649794303,48415.2766177886,46264.4235873111,48800.9285192888,47067.5005858836,47067.5,50964.1621028712,55183.4242066763,59751.9939801468,64698.7901880064,70055.1257449665,75854.9059245446,82134.8429773525,88934.6885173776,96297.4851472458,104269.83891524,112902.21432871,122249.253791289,132370.123485937,143328.887893225,155194.915315548,168043.316974196,181955.422458748,197019.294538355,213330.286593638,230991.646197698,250115.168666884,270821.904718247,293242.926713105,317520.158337007,343807.272967899,372270.666419123,403090.510214651,436461.892063708,472596.050753932)
for ($i = 0; $i -lt $vals104.Length; $i++) {
    $ws.Cells.Item(104, 10+$i).Value = $vals104[$i]
}

# row 111
$vals111 = @(2098.08481246737,2065.6808279352,1895.72422215299,1394.25182650581,1481.98232965245,1663.58565622714,1663.586,1771.26494299982,1885.91362171848,2007.98316629003,2137.95390715196,2276.33726509345,2423.67776364078,2580.55517169844,2747.58678487696,2925.42985448432,3114.78417373862,3316.39483137864,3531.05514350733,3759.60977520473,4002.95806419359,4262.05755963595,4537.92778998532,4831.65427472078,5144.39279574854,5477.37394527859,5831.90796807173,6209.38991711084,6611.30514298328,7039.23513857456,7494.86376207156,7979.98386276219,8496.50433570308,9046.45763301496,9632.00776136162,10255.4587970817)
for ($i = 0; $i -lt $vals111.Length; $i++) {
    $ws.Cells.Item(111, 10+$i).Value = $vals111[$i]
}

# row 112
$vals112 = @(62845.2603035521,70787.4299393758,69794.9233365629,74379.4615891143,79234.5722741823,72959.0251925868,72959.03,77494.8984180986,82312.7621191156,87430.1527711371,92865.6919873813,98639.1590881696,104771.563075666,111285.219083287,118203.829577986,125552.570610879,133358.183430087,141649.07178915,150455.405305117,159809.22924241,169744.581121969,180297.614579987,191506.730926954,203412.718885742,216058.903017212,229491.301373451,243758.79295233,258913.295562717,275009.954747587,292107.344452498,310267.680169636,329557.045333046,350045.63178887,371807.995215638,394923.32642406,419475.739523539)
for ($i = 0; $i -lt $vals112.Length; $i++) {
    $ws.Cells.Item(112, 10+$i).Value = $vals112[$i]
}

# --- Row 100 special handling ---
# J100:O100 -> new numeric values
$row100_first6 = @(9042.9928157095,1042.29264610104,8736.19277549749,10071.9445751427,6816.85776519906,24371.1743714232)
for ($i = 0; $i -lt $row100_first6.Length; $i++) {
    $ws.Cells.Item(100, 10+$i).Value = $row100_first6[$i]
}
# P100 -> rounded numeric value
$ws.Cells.Item(100, 16).Value = 24371.17
# Q100:AS100 -> text "inf"
for ($col = 17; $col -le 45; $col++) {
    $ws.Cells.Item(100, $col).Value = "inf"
}
